$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193, pushing existing rows 193:230 down to 194:231
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new data record
$ws.Cells.Item(193, 1).Value  = 10
$ws.Cells.Item(193, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(193, 3).Value  = "La Araucanía"
$ws.Cells.Item(193, 4).Value  = 44694
$ws.Cells.Item(193, 5).Value  = 9
$ws.Cells.Item(193, 6).Value  = "Fruta"
$ws.Cells.Item(193, 7).Value  = 100102
$ws.Cells.Item(193, 8).Value  = "Cítricos"
$ws.Cells.Item(193, 9).Value  = 100102006
$ws.Cells.Item(193, 10).Value = "Pomelo"
$ws.Cells.Item(193, 11).Value = "Start Ruby"
$ws.Cells.Item(193, 12).Value = "Primera"
$ws.Cells.Item(193, 13).Value = 25
$ws.Cells.Item(193, 14).Value = 13000
$ws.Cells.Item(193, 15).Value = 13000
$ws.Cells.Item(193, 16).Value = 13000
$ws.Cells.Item(193, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(193, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(193, 19).Value = 867
$ws.Cells.Item(193, 20).Value = 15

# Make sure the date cell keeps the date/time number format used by the rest of column D
$ws.Cells.Item(193, 4).NumberFormat = $ws.Cells.Item(194, 4).NumberFormat
